$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.199.62'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.655.88'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5312'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2626'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.41'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07829'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.523'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '1.660.09'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '1.883.09'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5498'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '0.0₅8168'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '26.173.11'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.601'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.39'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.09'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.027'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.009'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.64'
$ws.Range('E25').Value = '  +2.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1219'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.201'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.98'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('E29').Value = '  +4.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05776'
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.277'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.558'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.271'
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.593'
$ws.Range('E34').Value = '  +3.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.814'
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9531'
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.425'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5771'
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01602'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.815'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8523'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.046.63'
$ws.Range('E42').Value = '  +3.51%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.007'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.95'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('D45').Value = '1.795.23'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.85'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4365'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.881'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05156'
$ws.Range('E51').Value = '  +0.08%  '
